$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 'FAPs'
$ws.Cells.Item(2, 2).Value = 'Dlk2'
$ws.Cells.Item(2, 3).Value = 'Notch1'
$ws.Cells.Item(2, 4).Value = 'ECs'
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.426906
$ws.Cells.Item(2, 8).Value = 1.280718
$ws.Cells.Item(2, 9).Value = 0.7206774171221364
$ws.Cells.Item(2, 10).Value = 0.7206774171221364
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 58.95713633333333
$ws.Cells.Item(2, 14).Value = 176.871409
$ws.Cells.Item(2, 15).Value = 0.4863146960083892
$ws.Cells.Item(2, 16).Value = 0.4863146960083893
$ws.Cells.Item(2, 17).Value = 25.169155243518
$ws.Cells.Item(2, 18).Value = 226.522397191662
$ws.Cells.Item(2, 19).Value = 0.3504760190278629
$ws.Cells.Item(2, 20).Value = 0.3504760190278629

$ws.Cells.Item(3, 1).Value = 'FAPs'
$ws.Cells.Item(3, 2).Value = 'Dlk2'
$ws.Cells.Item(3, 3).Value = 'Notch1'
$ws.Cells.Item(3, 4).Value = 'FAPs'
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.426906
$ws.Cells.Item(3, 8).Value = 1.280718
$ws.Cells.Item(3, 9).Value = 0.7206774171221364
$ws.Cells.Item(3, 10).Value = 0.7206774171221364
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 8.99153
$ws.Cells.Item(3, 14).Value = 26.97459
$ws.Cells.Item(3, 15).Value = 0.07416766570679004
$ws.Cells.Item(3, 16).Value = 0.07416766570679005
$ws.Cells.Item(3, 17).Value = 3.83853810618
$ws.Cells.Item(3, 18).Value = 34.54684295562
$ws.Cells.Item(3, 19).Value = 0.0534509617555475
$ws.Cells.Item(3, 20).Value = 0.0534509617555475

$ws.Cells.Item(4, 1).Value = 'FAPs'
$ws.Cells.Item(4, 2).Value = 'Dlk2'
$ws.Cells.Item(4, 3).Value = 'Notch1'
$ws.Cells.Item(4, 4).Value = 'MuSCs'
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.426906
$ws.Cells.Item(4, 8).Value = 1.280718
$ws.Cells.Item(4, 9).Value = 0.7206774171221364
$ws.Cells.Item(4, 10).Value = 0.7206774171221364
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 42.51661933333333
$ws.Cells.Item(4, 14).Value = 127.549858
$ws.Cells.Item(4, 15).Value = 0.3507032073181665
$ws.Cells.Item(4, 16).Value = 0.3507032073181665
$ws.Cells.Item(4, 17).Value = 18.150599893116
$ws.Cells.Item(4, 18).Value = 163.355399038044
$ws.Cells.Item(4, 19).Value = 0.2527438816265054
$ws.Cells.Item(4, 20).Value = 0.2527438816265054

$ws.Cells.Item(5, 1).Value = 'FAPs'
$ws.Cells.Item(5, 2).Value = 'Dlk2'
$ws.Cells.Item(5, 3).Value = 'Notch1'
$ws.Cells.Item(5, 4).Value = 'Resolving-Mac'
$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 6).Value = 0.6666666666666666
$ws.Cells.Item(5, 7).Value = 0.426906
$ws.Cells.Item(5, 8).Value = 1.280718
$ws.Cells.Item(5, 9).Value = 0.7206774171221364
$ws.Cells.Item(5, 10).Value = 0.7206774171221364
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 10.76719366666667
$ws.Cells.Item(5, 14).Value = 32.301581
$ws.Cells.Item(5, 15).Value = 0.0888144309666542
$ws.Cells.Item(5, 16).Value = 0.08881443096665421
$ws.Cells.Item(5, 17).Value = 4.596579579461999
$ws.Cells.Item(5, 18).Value = 41.369216215158
$ws.Cells.Item(5, 19).Value = 0.06400655471222064
$ws.Cells.Item(5, 20).Value = 0.06400655471222065

$ws.Cells.Item(6, 1).Value = 'MuSCs'
$ws.Cells.Item(6, 2).Value = 'Dlk2'
$ws.Cells.Item(6, 3).Value = 'Notch1'
$ws.Cells.Item(6, 4).Value = 'ECs'
$ws.Cells.Item(6, 5).Value = 2
$ws.Cells.Item(6, 6).Value = 0.6666666666666666
$ws.Cells.Item(6, 7).Value = 0.1459986666666667
$ws.Cells.Item(6, 8).Value = 0.4379960000000001
$ws.Cells.Item(6, 9).Value = 0.2464662993647526
$ws.Cells.Item(6, 10).Value = 0.2464662993647526
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 58.95713633333333
$ws.Cells.Item(6, 14).Value = 176.871409
$ws.Cells.Item(6, 15).Value = 0.4863146960083892
$ws.Cells.Item(6, 16).Value = 0.4863146960083893
$ws.Cells.Item(6, 17).Value = 8.607663295151557
$ws.Cells.Item(6, 18).Value = 77.46896965636401
$ws.Cells.Item(6, 19).Value = 0.1198601834518823
$ws.Cells.Item(6, 20).Value = 0.1198601834518823

$ws.Cells.Item(7, 1).Value = 'MuSCs'
$ws.Cells.Item(7, 2).Value = 'Dlk2'
$ws.Cells.Item(7, 3).Value = 'Notch1'
$ws.Cells.Item(7, 4).Value = 'FAPs'
$ws.Cells.Item(7, 5).Value = 2
$ws.Cells.Item(7, 6).Value = 0.6666666666666666
$ws.Cells.Item(7, 7).Value = 0.1459986666666667
$ws.Cells.Item(7, 8).Value = 0.4379960000000001
$ws.Cells.Item(7, 9).Value = 0.2464662993647526
$ws.Cells.Item(7, 10).Value = 0.2464662993647526
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 8.99153
$ws.Cells.Item(7, 14).Value = 26.97459
$ws.Cells.Item(7, 15).Value = 0.07416766570679004
$ws.Cells.Item(7, 16).Value = 0.07416766570679005
$ws.Cells.Item(7, 17).Value = 1.312751391293334
$ws.Cells.Item(7, 18).Value = 11.81476252164
$ws.Cells.Item(7, 19).Value = 0.01827983009927461
$ws.Cells.Item(7, 20).Value = 0.01827983009927462

$ws.Cells.Item(8, 1).Value = 'MuSCs'
$ws.Cells.Item(8, 2).Value = 'Dlk2'
$ws.Cells.Item(8, 3).Value = 'Notch1'
$ws.Cells.Item(8, 4).Value = 'MuSCs'
$ws.Cells.Item(8, 5).Value = 2
$ws.Cells.Item(8, 6).Value = 0.6666666666666666
$ws.Cells.Item(8, 7).Value = 0.1459986666666667
$ws.Cells.Item(8, 8).Value = 0.4379960000000001
$ws.Cells.Item(8, 9).Value = 0.2464662993647526
$ws.Cells.Item(8, 10).Value = 0.2464662993647526
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 42.51661933333333
$ws.Cells.Item(8, 14).Value = 127.549858
$ws.Cells.Item(8, 15).Value = 0.3507032073181665
$ws.Cells.Item(8, 16).Value = 0.3507032073181665
$ws.Cells.Item(8, 17).Value = 6.207369733840889
$ws.Cells.Item(8, 18).Value = 55.86632760456801
$ws.Cells.Item(8, 19).Value = 0.08643652168305813
$ws.Cells.Item(8, 20).Value = 0.08643652168305814

$ws.Cells.Item(9, 1).Value = 'MuSCs'
$ws.Cells.Item(9, 2).Value = 'Dlk2'
$ws.Cells.Item(9, 3).Value = 'Notch1'
$ws.Cells.Item(9, 4).Value = 'Resolving-Mac'
$ws.Cells.Item(9, 5).Value = 2
$ws.Cells.Item(9, 6).Value = 0.6666666666666666
$ws.Cells.Item(9, 7).Value = 0.1459986666666667
$ws.Cells.Item(9, 8).Value = 0.4379960000000001
$ws.Cells.Item(9, 9).Value = 0.2464662993647526
$ws.Cells.Item(9, 10).Value = 0.2464662993647526
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 10.76719366666667
$ws.Cells.Item(9, 14).Value = 32.301581
$ws.Cells.Item(9, 15).Value = 0.0888144309666542
$ws.Cells.Item(9, 16).Value = 0.08881443096665421
$ws.Cells.Item(9, 17).Value = 1.571995919075111
$ws.Cells.Item(9, 18).Value = 14.147963271676
$ws.Cells.Item(9, 19).Value = 0.02188976413053755
$ws.Cells.Item(9, 20).Value = 0.02188976413053756

$ws.Cells.Item(10, 1).Value = 'Resolving-Mac'
$ws.Cells.Item(10, 2).Value = 'Dlk2'
$ws.Cells.Item(10, 3).Value = 'Notch1'
$ws.Cells.Item(10, 4).Value = 'ECs'
$ws.Cells.Item(10, 5).Value = 2
$ws.Cells.Item(10, 6).Value = 0.6666666666666666
$ws.Cells.Item(10, 7).Value = 0.019463
$ws.Cells.Item(10, 8).Value = 0.058389
$ws.Cells.Item(10, 9).Value = 0.03285628351311094
$ws.Cells.Item(10, 10).Value = 0.03285628351311094
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 58.95713633333333
$ws.Cells.Item(10, 14).Value = 176.871409
$ws.Cells.Item(10, 15).Value = 0.4863146960083892
$ws.Cells.Item(10, 16).Value = 0.4863146960083893
$ws.Cells.Item(10, 17).Value = 1.147482744455667
$ws.Cells.Item(10, 18).Value = 10.327344700101
$ws.Cells.Item(10, 19).Value = 0.015978493528644
$ws.Cells.Item(10, 20).Value = 0.015978493528644

$ws.Cells.Item(11, 1).Value = 'Resolving-Mac'
$ws.Cells.Item(11, 2).Value = 'Dlk2'
$ws.Cells.Item(11, 3).Value = 'Notch1'
$ws.Cells.Item(11, 4).Value = 'FAPs'
$ws.Cells.Item(11, 5).Value = 2
$ws.Cells.Item(11, 6).Value = 0.6666666666666666
$ws.Cells.Item(11, 7).Value = 0.019463
$ws.Cells.Item(11, 8).Value = 0.058389
$ws.Cells.Item(11, 9).Value = 0.03285628351311094
$ws.Cells.Item(11, 10).Value = 0.03285628351311094
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 8.99153
$ws.Cells.Item(11, 14).Value = 26.97459
$ws.Cells.Item(11, 15).Value = 0.07416766570679004
$ws.Cells.Item(11, 16).Value = 0.07416766570679005
$ws.Cells.Item(11, 17).Value = 0.17500214839
$ws.Cells.Item(11, 18).Value = 1.57501933551
$ws.Cells.Item(11, 19).Value = 0.002436873851967929
$ws.Cells.Item(11, 20).Value = 0.00243687385196793

$ws.Cells.Item(12, 1).Value = 'Resolving-Mac'
$ws.Cells.Item(12, 2).Value = 'Dlk2'
$ws.Cells.Item(12, 3).Value = 'Notch1'
$ws.Cells.Item(12, 4).Value = 'MuSCs'
$ws.Cells.Item(12, 5).Value = 2
$ws.Cells.Item(12, 6).Value = 0.6666666666666666
$ws.Cells.Item(12, 7).Value = 0.019463
$ws.Cells.Item(12, 8).Value = 0.058389
$ws.Cells.Item(12, 9).Value = 0.03285628351311094
$ws.Cells.Item(12, 10).Value = 0.03285628351311094
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 42.51661933333333
$ws.Cells.Item(12, 14).Value = 127.549858
$ws.Cells.Item(12, 15).Value = 0.3507032073181665
$ws.Cells.Item(12, 16).Value = 0.3507032073181665
$ws.Cells.Item(12, 17).Value = 0.8275009620846665
$ws.Cells.Item(12, 18).Value = 7.447508658762
$ws.Cells.Item(12, 19).Value = 0.011522804008603
$ws.Cells.Item(12, 20).Value = 0.011522804008603

$ws.Cells.Item(13, 1).Value = 'Resolving-Mac'
$ws.Cells.Item(13, 2).Value = 'Dlk2'
$ws.Cells.Item(13, 3).Value = 'Notch1'
$ws.Cells.Item(13, 4).Value = 'Resolving-Mac'
$ws.Cells.Item(13, 5).Value = 2
$ws.Cells.Item(13, 6).Value = 0.6666666666666666
$ws.Cells.Item(13, 7).Value = 0.019463
$ws.Cells.Item(13, 8).Value = 0.058389
$ws.Cells.Item(13, 9).Value = 0.03285628351311094
$ws.Cells.Item(13, 10).Value = 0.03285628351311094
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 10.76719366666667
$ws.Cells.Item(13, 14).Value = 32.301581
$ws.Cells.Item(13, 15).Value = 0.0888144309666542
$ws.Cells.Item(13, 16).Value = 0.08881443096665421
$ws.Cells.Item(13, 17).Value = 0.2095618903343333
$ws.Cells.Item(13, 18).Value = 1.886057013009
$ws.Cells.Item(13, 19).Value = 0.00291811212389601
$ws.Cells.Item(13, 20).Value = 0.002918112123896011

"done"